$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the QPU parameter values (row 8)
$ws.Range("B8").Value = 160
$ws.Range("C8").Value = 130
$ws.Range("D8").Value = 120

# Update the active cell selection to match the saved view state (F8)
$ws.Range("F8").Select()

$wb.Save()
